# Auto-generated Excel COM-interop script
# Applies numeric cell updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 204.6
$ws.Range("I2").Value = 125
$ws.Range("J2").Value = 257.66666
$ws.Range("K2").Value = 125
$ws.Range("L2").Value = 257.66666
$ws.Range("M2").Value = -12
$ws.Range("N2").Value = -483.66666
$ws.Range("H38").Value = 59059.883
$ws.Range("J38").Value = 900
$ws.Range("L38").Value = 2700
$ws.Range("N38").Value = -3444
$ws.Range("H58").Value = 2103.0667
$ws.Range("I58").Value = 1258.5714
$ws.Range("J58").Value = 2842
$ws.Range("K58").Value = 3775.7142
$ws.Range("L58").Value = 8526
$ws.Range("M58").Value = -3625.7142
$ws.Range("N58").Value = -8826
$ws.Range("H98").Value = 1228.9166
$ws.Range("I98").Value = 1124.1
$ws.Range("J98").Value = 1753
$ws.Range("K98").Value = 1124.1
$ws.Range("L98").Value = 1753
$ws.Range("M98").Value = 373.9000000000001
$ws.Range("N98").Value = -4749
$ws.Range("H122").Value = 1228.9166
$ws.Range("I122").Value = 1124.1
$ws.Range("J122").Value = 1753
$ws.Range("K122").Value = 3372.3
$ws.Range("L122").Value = 5259
$ws.Range("M122").Value = -922.2999999999997
$ws.Range("N122").Value = -10159
$ws.Range("H137").Value = 1186.9667
$ws.Range("I137").Value = 937.2174
$ws.Range("J137").Value = 2007.5714
$ws.Range("K137").Value = 2811.6522
$ws.Range("L137").Value = 6022.7142
$ws.Range("M137").Value = -261.6522
$ws.Range("N137").Value = -11122.7142
$ws.Range("H138").Value = 3194.4736
$ws.Range("J138").Value = 3113.8462
$ws.Range("L138").Value = 9341.5386
$ws.Range("N138").Value = -19621.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 412322.4
$ws.Range("I32").Value = 461338.56
$ws.Range("J32").Value = 20193.125
$ws.Range("K32").Value = 461338.56
$ws.Range("L32").Value = 20193.125
$ws.Range("M32").Value = -461051.56
$ws.Range("N32").Value = -20767.125
$ws.Range("H61").Value = 6946661.5
$ws.Range("I61").Value = 15152731
$ws.Range("J61").Value = 3064.0386
$ws.Range("K61").Value = 15152731
$ws.Range("L61").Value = 3064.0386
$ws.Range("M61").Value = -15152519
$ws.Range("N61").Value = -3488.0386
$ws.Range("H74").Value = 856.45
$ws.Range("I74").Value = 636.9167
$ws.Range("J74").Value = 950.5357
$ws.Range("K74").Value = 636.9167
$ws.Range("L74").Value = 950.5357
$ws.Range("M74").Value = 237.0833
$ws.Range("N74").Value = -2698.5357
$ws.Range("H77").Value = 856.45
$ws.Range("I77").Value = 636.9167
$ws.Range("J77").Value = 950.5357
$ws.Range("K77").Value = 3184.5835
$ws.Range("L77").Value = 4752.6785
$ws.Range("M77").Value = 1183.4165
$ws.Range("N77").Value = -13488.6785
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 5188.452
$ws.Range("I132").Value = 5100.0967
$ws.Range("J132").Value = 5437.4546
$ws.Range("K132").Value = 15300.2901
$ws.Range("L132").Value = 16312.3638
$ws.Range("M132").Value = -12770.2901
$ws.Range("N132").Value = -21372.3638
$ws.Range("H133").Value = 59420.332
$ws.Range("J133").Value = 59420.332
$ws.Range("L133").Value = 59420.332
$ws.Range("N133").Value = -64480.332
$ws.Range("H136").Value = 6946661.5
$ws.Range("I136").Value = 15152731
$ws.Range("J136").Value = 3064.0386
$ws.Range("K136").Value = 45458193
$ws.Range("L136").Value = 9192.1158
$ws.Range("M136").Value = -45455643
$ws.Range("N136").Value = -14292.1158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3242.6
$ws.Range("I86").Value = 1735.3334
$ws.Range("J86").Value = 5503.5
$ws.Range("K86").Value = 1735.3334
$ws.Range("L86").Value = 5503.5
$ws.Range("M86").Value = -612.3334
$ws.Range("N86").Value = -7749.5
$ws.Range("H89").Value = 3242.6
$ws.Range("I89").Value = 1735.3334
$ws.Range("J89").Value = 5503.5
$ws.Range("K89").Value = 8676.666999999999
$ws.Range("L89").Value = 27517.5
$ws.Range("M89").Value = -3060.666999999999
$ws.Range("N89").Value = -38749.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1502.415
$ws.Range("I31").Value = 1081.5217
$ws.Range("J31").Value = 1825.1
$ws.Range("K31").Value = 1081.5217
$ws.Range("L31").Value = 1825.1
$ws.Range("M31").Value = -786.5217
$ws.Range("N31").Value = -2415.1
$ws.Range("H34").Value = 1502.415
$ws.Range("I34").Value = 1081.5217
$ws.Range("J34").Value = 1825.1
$ws.Range("K34").Value = 1081.5217
$ws.Range("L34").Value = 1825.1
$ws.Range("M34").Value = -879.5217
$ws.Range("N34").Value = -2229.1
$ws.Range("H60").Value = 14013.333
$ws.Range("J60").Value = 17496.75
$ws.Range("L60").Value = 17496.75
$ws.Range("N60").Value = -18518.75
$ws.Range("H99").Value = 1801.7
$ws.Range("I99").Value = 1339
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1339
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 159
$ws.Range("N99").Value = -4996
$ws.Range("H126").Value = 1801.7
$ws.Range("I126").Value = 1339
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4017
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1547
$ws.Range("N126").Value = -10940
$ws.Range("H129").Value = 32449.5
$ws.Range("I129").Value = 14900
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 14900
$ws.Range("L129").Value = 49999
$ws.Range("M129").Value = -9900
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 12823404
$ws.Range("I132").Value = 3544.75
$ws.Range("J132").Value = 18521118
$ws.Range("K132").Value = 10634.25
$ws.Range("L132").Value = 55563354
$ws.Range("M132").Value = -8104.25
$ws.Range("N132").Value = -55568414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2875
$ws.Range("H131").Value = 1057.6957
$ws.Range("I131").Value = 847.5714
$ws.Range("J131").Value = 1149.625
$ws.Range("K131").Value = 2542.7142
$ws.Range("L131").Value = 3448.875
$ws.Range("M131").Value = 2497.2858
$ws.Range("N131").Value = -13528.875
$ws.Range("H137").Value = 13926.1
$ws.Range("I137").Value = 24565.6
$ws.Range("J137").Value = 3286.6
$ws.Range("K137").Value = 73696.79999999999
$ws.Range("L137").Value = 9859.799999999999
$ws.Range("M137").Value = -68596.79999999999
$ws.Range("N137").Value = -20059.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 31333.334
$ws.Range("J15").Value = 31333.334
$ws.Range("L15").Value = 31333.334
$ws.Range("N15").Value = -31909.334
$ws.Range("H81").Value = 31333.334
$ws.Range("J81").Value = 31333.334
$ws.Range("L81").Value = 31333.334
$ws.Range("N81").Value = -33329.334
$ws.Range("H84").Value = 31333.334
$ws.Range("J84").Value = 31333.334
$ws.Range("L84").Value = 94000.00199999999
$ws.Range("N84").Value = -103984.002
$ws.Range("H122").Value = 2151.3
$ws.Range("I122").Value = 1899.25
$ws.Range("K122").Value = 5697.75
$ws.Range("M122").Value = -3247.75
$ws.Range("H132").Value = 4999.3335
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 166669150
$ws.Range("I7").Value = 250002260
$ws.Range("J7").Value = 2952.5
$ws.Range("K7").Value = 250002260
$ws.Range("L7").Value = 2952.5
$ws.Range("M7").Value = -250002148
$ws.Range("N7").Value = -3176.5
$ws.Range("H40").Value = 29415346
$ws.Range("I40").Value = 52633616
$ws.Range("J40").Value = 5540
$ws.Range("K40").Value = 52633616
$ws.Range("L40").Value = 5540
$ws.Range("M40").Value = -52633480
$ws.Range("N40").Value = -5812
$ws.Range("H61").Value = 4426.5625
$ws.Range("I61").Value = 4186.25
$ws.Range("J61").Value = 4666.875
$ws.Range("K61").Value = 4186.25
$ws.Range("L61").Value = 4666.875
$ws.Range("M61").Value = -3984.25
$ws.Range("N61").Value = -5070.875
$ws.Range("H113").Value = 4426.5625
$ws.Range("I113").Value = 4186.25
$ws.Range("J113").Value = 4666.875
$ws.Range("K113").Value = 4186.25
$ws.Range("L113").Value = 4666.875
$ws.Range("M113").Value = -2016.25
$ws.Range("N113").Value = -9006.875
$ws.Range("H122").Value = 2064
$ws.Range("I122").Value = 1876.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5630.4
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3180.4
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 166669150
$ws.Range("I126").Value = 250002260
$ws.Range("J126").Value = 2952.5
$ws.Range("K126").Value = 750006780
$ws.Range("L126").Value = 8857.5
$ws.Range("M126").Value = -750004310
$ws.Range("N126").Value = -13797.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1564.875
$ws.Range("I126").Value = 1001
$ws.Range("J126").Value = 2128.75
$ws.Range("K126").Value = 3003
$ws.Range("L126").Value = 6386.25
$ws.Range("M126").Value = -533
$ws.Range("N126").Value = -11326.25
$ws.Range("H136").Value = 1819.961
$ws.Range("I136").Value = 1533.9259
$ws.Range("K136").Value = 4601.7777
$ws.Range("M136").Value = -2051.7777
